$d = $word.ActiveDocument

# Update the date paragraph
$d.Content.Find.Execute("2025-04-09 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-04-10 Thursday", 2) | Out-Null

# Update all 100 table cells (20 rows x 5 columns), in row-major order
$newValues = @(
    "76+2=",
    "87-30=",
    "93+4=",
    "47-7=",
    "71-20=",
    "63+17=",
    "5+31=",
    "61+30=",
    "45-5=",
    "53+37=",
    "19-5=",
    "28-9=",
    "14-8=",
    "56+33=",
    "82+2=",
    "56-48=",
    "53-0=",
    "76-41=",
    "88-18=",
    "36-22=",
    "6+43=",
    "91-49=",
    "54-43=",
    "22+16=",
    "87-74=",
    "85-18=",
    "8-4=",
    "94-22=",
    "70+2=",
    "86-66=",
    "65+24=",
    "53-4=",
    "50-17=",
    "92+4=",
    "59+22=",
    "37+22=",
    "85-69=",
    "95-13=",
    "52-42=",
    "90-17=",
    "46+24=",
    "57+31=",
    "76-48=",
    "2+84=",
    "93-84=",
    "5+21=",
    "16+75=",
    "33-27=",
    "63-6=",
    "18+9=",
    "98-92=",
    "26+46=",
    "83-42=",
    "41-19=",
    "44-40=",
    "7+29=",
    "58+14=",
    "64-49=",
    "88-58=",
    "7+76=",
    "49+6=",
    "76-66=",
    "28+56=",
    "48-48=",
    "10-9=",
    "54+19=",
    "68+1=",
    "8+66=",
    "46-22=",
    "99-95=",
    "54-3=",
    "26+3=",
    "79-61=",
    "50+40=",
    "58+13=",
    "47-19=",
    "67+4=",
    "6+38=",
    "24+5=",
    "35-25=",
    "79-4=",
    "25+67=",
    "15+58=",
    "25-17=",
    "22+35=",
    "62-6=",
    "78+2=",
    "66-60=",
    "16+61=",
    "30-29=",
    "56-37=",
    "54+13=",
    "5+89=",
    "96-55=",
    "88-48=",
    "9+28=",
    "36-16=",
    "19-10=",
    "16+46=",
    "46+5="
)

$t = $d.Tables(1)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "Updated $idx cells"